$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date serials (Excel 1900 date system) for 2004Q4..2024Q4 (i.e. YYYY-12-31)
$serials = @(38352,38717,39082,39447,39813,40178,40543,40908,41274,41639,42004,42369,42735,43100,43465,43830,44196,44561,44926,45291,45657)

for ($i = 0; $i -lt $serials.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value2 = $serials[$i]
    $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
